# Daily update at 8 AM UTC
# Appends the new day's row (row 65) to the Wins Over Time tracker and
# shifts the "last row" date-only formatting down from row 64 to row 65.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64 was previously the last row (date-only format). It is no longer the
# last row, so it reverts to the regular date-time number format used by all
# the other data rows.
$ws.Range("A64").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New last row with the next day's data.
$ws.Range("A65").Value = 45805
$ws.Range("B65").Value = 275
$ws.Range("C65").Value = 280
$ws.Range("D65").Value = 279

# The new last row picks up the date-only number format previously on A64.
$ws.Range("A65").NumberFormat = "YYYY-MM-DD"
